# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the latest scraped data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row number => new F value
$exhibitionUpdates = @{
    2  = 1892
    3  = 508
    6  = 2672
    10 = 1563
    11 = 545
    13 = 339
    17 = 2
    22 = 207
    23 = 67
    24 = 1716
    25 = 38
    26 = 417
    27 = 53
    30 = 307
    31 = 435
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" - row number => new F value
$allTypesUpdates = @{
    2  = 1892
    4  = 508
    7  = 2672
    11 = 1563
    12 = 545
    14 = 339
    18 = 2
    23 = 207
    24 = 67
    25 = 1716
    26 = 38
    27 = 417
    28 = 53
    31 = 307
    32 = 435
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
